# Apply the edits described by the diff:
#  - Rename the two "round" header pairs (Ngay Duyet/Tu choi + Trang thai)
#    from numbered ("1" / "2") to named ("Lan dau" / "Lan cuoi") labels.
#  - Move the active cell selection on the sheet from E9 to D6.
#  - Slightly widen columns D and I.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text (shared strings) -----------------------------------
$ws.Range("G1").Value = "Ngày Duyệt/Từ chối Lần đầu"
$ws.Range("H1").Value = "Trạng thái Lần đầu"
$ws.Range("I1").Value = "Ngày Duyệt/Từ chối Lần cuối"
$ws.Range("J1").Value = "Trạng thái Lần cuối"

# --- Resize columns D (4) and I (9) -----------------------------------------
# Target stored widths (OOXML <col width="...">) are 12.2814814814815 and
# 24.7777777777778 respectively; this runtime quantizes ColumnWidth to whole
# pixels (character-width units of 1/7), so we pick the ColumnWidth value
# that lands on the closest achievable stored width.
$ws.Columns.Item(4).ColumnWidth = 11.571428571428571
$ws.Columns.Item(9).ColumnWidth = 24.0

# --- Move the selected / active cell ----------------------------------------
$null = $ws.Range("D6").Select()
